# Bold the text in the last row of the benefit/feature tables on several
# slides (Slide 4, 8, 12, 17, 21). Each of these slides has a single table
# placeholder (Shape 3) and the final data row of that table should have
# its run text made bold, matching the header row's bold treatment.

$p = $ppt.ActivePresentation

$slideNumbers = @(4, 8, 12, 17, 21)

foreach ($slideNumber in $slideNumbers) {
    $slide = $p.Slides.Item($slideNumber)

    # Find the shape that owns the table on this slide.
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            $lastRow = $table.Rows.Count
            $colCount = $table.Columns.Count

            for ($col = 1; $col -le $colCount; $col++) {
                $cell = $table.Cell($lastRow, $col)
                $cell.Shape.TextFrame.TextRange.Font.Bold = $true
            }
        }
    }
}
